$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$b2 = 1.063418937352623 / 10000000
$c2 = 0.3375848360084654
$d2 = 157.8057217802531
$e2 = 1.594453305621061 * 1000000000000000000
$f2 = 1
$g2 = 1.594453305621061 * 1000000000000000000

$ws.Range("B2").Value = $b2
$ws.Range("C2").Value = $c2
$ws.Range("D2").Value = $d2
$ws.Range("E2").Value = $e2
$ws.Range("F2").Value = $f2
$ws.Range("G2").Value = $g2
